# Replace the "Test" cycle (rows 5-6) with a new "Esame 2" cycle: new
# timestamps, new flow readings, and Maturation phase/status on row 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 1).Value = "Esame 2"
$ws.Cells.Item(5, 2).Value = "2018-09-10 06:33:36.694791"
$ws.Cells.Item(5, 3).Value = "2018-09-10 06:33:57.893403"
$ws.Cells.Item(5, 4).Value = 20
$ws.Cells.Item(5, 5).Value = 28
$ws.Cells.Item(5, 6).Value = 20
$ws.Cells.Item(5, 8).Value = "Casting: OK"

# Row 6
$ws.Cells.Item(6, 1).Value = "Esame 2"
$ws.Cells.Item(6, 2).Value = "2018-09-10 06:33:36.694791"
$ws.Cells.Item(6, 3).Value = "2018-09-10 06:34:28.174998"
$ws.Cells.Item(6, 4).Value = "'20.0"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'40.0"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 6).Value = "'20.0"
$ws.Cells.Item(6, 6).Style = "Normal"
$ws.Cells.Item(6, 7).Value = "Maturation"
$ws.Cells.Item(6, 8).Value = "Maturation: OK"
